$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 5137677138
$ws.Range("B2").Value = 0
